$wb = $excel.ActiveWorkbook

# Sheet "展览" (index 1) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 15
$ws1.Range("F5").Value = 6679
$ws1.Range("F6").Value = 5473
$ws1.Range("F12").Value = 121
$ws1.Range("F13").Value = 44

# Sheet "全部类型" (index 4) updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 15
$ws4.Range("F5").Value = 6679
$ws4.Range("F6").Value = 5473
$ws4.Range("F14").Value = 121
$ws4.Range("F15").Value = 44
